{"js": "// Office.js (Word JavaScript API) script\n// Applies the text replacements described by the diff.\n\nconst replacements = [\n  [\n    \"Play Ghost Glyph Slot for Free: Cluster Pays with Tumble & Urn Features\",\n    \"Play Ghost Glyph for Free - Spooky Slot Fun\"\n  ],\n  [\n    \"Tumble mechanic allows for avalanche wins\",\n    \"Cluster Pays mechanic for exciting gameplay\"\n  ],\n  [\n    \"Cluster Pays offers unique gameplay\",\n    \"Tumble mechanic adds to the thrill of winning\"\n  ],\n  [\n    \"Innovative urn features add excitement\",\n    \"Ghost Glyphs as wild symbols and added to urns\"\n  ],\n  [\n    \"Massive Wild feature with 22 levels\",\n    \"Various urn features with unique benefits\"\n  ],\n  [\n    \"No progressive jackpot\",\n    \"Limited number of categories for symbols\"\n  ],\n  [\n    \"Bonus game with free spins can be infrequent\",\n    \"Massive Wild feature requires specific conditions to activate\"\n  ],\n  [\n    \"Get ready to play Ghost Glyph online slot for free! This game offers a cluster pays mechanism, Tumble mechanic, urn features, Massive Wild, a spooky yet playful theme, a 96.00% RTP, and a betting range starting at 20c.\",\n    \"Play Ghost Glyph for free and experience thrilling gameplay with ghostly symbols and exciting features.\"\n  ]\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) script\n# Applies the text replacements described by the diff.\n\n$d = $word.ActiveDocument\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\nfunction Replace-Text($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($findText, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $replaceText, $wdReplaceAll)\n}\n\nReplace-Text \"Play Ghost Glyph Slot for Free: Cluster Pays with Tumble & Urn Features\" \"Play Ghost Glyph for Free - Spooky Slot Fun\"\nReplace-Text \"Tumble mechanic allows for avalanche wins\" \"Cluster Pays mechanic for exciting gameplay\"\nReplace-Text \"Cluster Pays offers unique gameplay\" \"Tumble mechanic adds to the thrill of winning\"\nReplace-Text \"Innovative urn features add excitement\" \"Ghost Glyphs as wild symbols and added to urns\"\nReplace-Text \"Massive Wild feature with 22 levels\" \"Various urn features with unique benefits\"\nReplace-Text \"No progressive jackpot\" \"Limited number of categories for symbols\"\nReplace-Text \"Bonus game with free spins can be infrequent\" \"Massive Wild feature requires specific conditions to activate\"\nReplace-Text \"Get ready to play Ghost Glyph online slot for free! This game offers a cluster pays mechanism, Tumble mechanic, urn features, Massive Wild, a spooky yet playful theme, a 96.00% RTP, and a betting range starting at 20c.\" \"Play Ghost Glyph for free and experience thrilling gameplay with ghostly symbols and exciting features.\"\n"}
